$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 'hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G3").Value = 'Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range("G4").Value = 'Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G5").Value = 'Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G6").Value = 'gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G7").Value = 'Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G8").Value = 'backup@backdoor.com, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G9").Value = 'gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G10").Value = 'drsarahmw@gmail.com, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, heba@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G12").Value = 'salma.elgendy.std@med.asu.edu.eg, System'
$ws.Range("G13").Value = 'wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg'
$ws.Range("G14").Value = 'Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G20").Value = 'user@user.com, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G24").Value = 'maryam.ashraf@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G25").Value = 'marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G27").Value = 'neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G28").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G29").Value = 'Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
$ws.Range("G30").Value = 'Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G31").Value = 'Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg'
$ws.Range("G32").Value = 'Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G33").Value = 'Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G34").Value = 'gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range("G35").Value = 'Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G36").Value = 'backup@backdoor.com, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Range("G37").Value = 'gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G38").Value = 'drsarahmw@gmail.com, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, heba@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg'
$ws.Range("G40").Value = 'salma.elgendy.std@med.asu.edu.eg, System'
$ws.Range("G41").Value = 'wessam.atef@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg'
$ws.Range("G42").Value = 'Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G48").Value = 'user@user.com, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G52").Value = 'maryam.ashraf@med.asu.edu.eg, marina_atef@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G53").Value = 'marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, eman.samir@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G55").Value = 'neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Eman_mohamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G56").Value = 'neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg'
$ws.Range("G57").Value = 'Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg'
